$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list (prices + 1h volume change) to the latest snapshot.
# Note: some Price values are plain decimals (e.g. "606.52"); assigning
# those via .Value would make Excel auto-convert them to numbers, but the
# source data stores every Price/Volume cell as text. A leading quote
# ('' below, i.e. a literal single quote doubled for PowerShell's
# single-quoted string escaping) forces Excel to keep them as text, the
# same way typing '606.52 into a cell does.

$ws.Range('D2').Value = '71.632.51'
$ws.Range('E2').Value = '  +4.15%  '
$ws.Range('D3').Value = '2.626.88'
$ws.Range('E3').Value = '  +4.45%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''606.52'
$ws.Range('E5').Value = '  +2.31%  '
$ws.Range('D6').Value = '''179.22'
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('D9').Value = '2.625.99'
$ws.Range('E9').Value = '  +4.44%  '
$ws.Range('D10').Value = '''0.168'
$ws.Range('E10').Value = '  +14.15%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('E12').Value = '  +3.09%  '
$ws.Range('D13').Value = '''5.04'
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('D14').Value = '3.134.52'
$ws.Range('E14').Value = '  +4.70%  '
$ws.Range('D15').Value = '''0.0000186'
$ws.Range('E15').Value = '  +9.02%  '
$ws.Range('D16').Value = '''26.55'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '71.563.67'
$ws.Range('E17').Value = '  +4.15%  '
$ws.Range('D18').Value = '2.637.81'
$ws.Range('E18').Value = '  +5.06%  '
$ws.Range('D19').Value = '''382.94'
$ws.Range('E19').Value = '  +6.20%  '
$ws.Range('D20').Value = '''7.96'
$ws.Range('E20').Value = '  +6.40%  '
$ws.Range('D21').Value = '''11.45'
$ws.Range('E21').Value = '  +4.54%  '
$ws.Range('D22').Value = '''4.16'
$ws.Range('E22').Value = '  +2.08%  '
$ws.Range('E23').Value = '  +18.37%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''72.57'
$ws.Range('E24').Value = '  +2.85%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '''4.46'
$ws.Range('E25').Value = '  +6.95%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '''9.97'
$ws.Range('E27').Value = '  +10.74%  '
$ws.Range('D28').Value = '2.763.09'
$ws.Range('E28').Value = '  +4.39%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0961'
$ws.Range('E30').Value = '  +9.53%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '''544.25'
$ws.Range('E31').Value = '  +6.89%  '
$ws.Range('D32').Value = '''8.04'
$ws.Range('E32').Value = '  +3.91%  '
$ws.Range('D33').Value = '''1.33'
$ws.Range('E33').Value = '  +8.39%  '
$ws.Range('E34').Value = '  +3.51%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '''166.22'
$ws.Range('E36').Value = '  +3.72%  '
$ws.Range('D37').Value = '''19.20'
$ws.Range('E37').Value = '  +3.17%  '
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').Value = '''19.16'
$ws.Range('E39').Value = '  +2.73%  '
$ws.Range('E40').Value = '  +6.75%  '
$ws.Range('D41').Value = '''1.86'
$ws.Range('E41').Value = '  +8.24%  '
$ws.Range('E42').Value = '  +11.96%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('E44').Value = '  +5.40%  '
$ws.Range('D45').Value = '''0.332'
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D47').Value = '''151.07'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('D49').Value = '''0.535'
$ws.Range('E49').Value = '  +4.33%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = '''1.69'
$ws.Range('E50').Value = '  +7.88%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0264'
$ws.Range('E51').Value = '  +5.90%  '
